# Sync attendance_reports: swap the order of "dnasr281@gmail.com" and
# "System" inside the G column "Recorded By" text for a specific subset
# of session rows (matching the upstream commit's canonical diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$targetRows = @(
    2,3,4,5,6,7,
    16,17,
    22,23,
    37,38,
    43,44,
    58,59,
    64,65,
    79,80,
    85,86,87,88,89,90,
    99,100,
    105,106,107,108,109,110,
    119,120,
    125,126,127,128,129,130,
    139,140,
    145,146,147,148,149,150,
    159,160,
    165,166,167,168,169,170,
    179,180,
    185,186,
    200,201,
    206,207,
    221,222,
    227,228,
    242,243
)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
